$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 "土地" (land): clean up stray spaces/hyphens in existing text and
# append the new metadata columns (I:O) used by the property scraper output.
# ---------------------------------------------------------------------------
$landSheet = $wb.Worksheets.Item(1)

$landSheet.Range("B2").Value = "新北市永和區永福段08230000地號"
$landSheet.Range("F2").Value = "80年06月13曰"
$landSheet.Range("G2").Value = "第一次登記"

$landSheet.Range("I1").Value = "property_category"
$landSheet.Range("J1").Value = "category"
$landSheet.Range("K1").Value = "date"
$landSheet.Range("L1").Value = "legislator_name"
$landSheet.Range("M1").Value = "legislator_id"
$landSheet.Range("N1").Value = "source_file"
$landSheet.Range("O1").Value = "index"

$landSheet.Range("I2").Value = "land"
$landSheet.Range("J2").Value = "normal"
# Force the date-shaped value to stay literal text instead of being
# auto-parsed into a date serial number by the COM layer.
$landSheet.Range("K2").NumberFormat = "@"
$landSheet.Range("K2").Value = "2013-12-30"
$landSheet.Range("L2").Value = "林淑芬"
$landSheet.Range("M2").Value = 1337
$landSheet.Range("N2").Value = "tmp63cf1"
$landSheet.Range("O2").Value = 14

# ---------------------------------------------------------------------------
# Sheet 2 "建物" (building): same text clean-up, no structural changes.
# ---------------------------------------------------------------------------
$buildingSheet = $wb.Worksheets.Item(2)

$buildingSheet.Range("B2").Value = "新北市永和區永福段02017000建號"
$buildingSheet.Range("F2").Value = "80年06月13日"
$buildingSheet.Range("G2").Value = "第一次登記"

# ---------------------------------------------------------------------------
# Sheet 3 "汽車" (car): same text clean-up.
# ---------------------------------------------------------------------------
$carSheet = $wb.Worksheets.Item(3)

$carSheet.Range("B2").Value = "toyotarav4rod"
$carSheet.Range("E2").Value = "100年03月01曰"

# ---------------------------------------------------------------------------
# Sheet 4 "現金" (cash): unchanged.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Sheet 5 "存款" (deposit): text clean-up + F4 becomes a real number.
# ---------------------------------------------------------------------------
$depositSheet = $wb.Worksheets.Item(5)

$depositSheet.Range("B3").Value = "遠東國際商業銀行永和分行"
$depositSheet.Range("B4").Value = "中國信託商業銀行美國中信銀子行"
$depositSheet.Range("F4").Value = 1420000
